# Applies the committed changes to "web pitch deck.pptx":
#  1. Refresh the cached `datetimeFigureOut` field text (9/11/2024 -> 9/22/2024)
#     on every slide layout's and the slide master's Date Placeholder shape.
#  2. Trim the feature-summary textbox on slide 4 and resize it to its
#     auto-fit height for the shorter text.
#  3. Trim the "By utilizing Diet Tracker" textbox on slide 8 and resize it
#     to its auto-fit height for the shorter text.

$p = $ppt.ActivePresentation

function Get-ShapeByName($container, $name) {
    for ($k = 1; $k -le $container.Shapes.Count; $k++) {
        $s = $container.Shapes.Item($k)
        if ($s.Name -eq $name) {
            return $s
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Date placeholder text on the slide master + every custom layout
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container, $newDateText) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = $newDateText
            }
        }
    }
}

$master = $p.Designs.Item(1).SlideMaster
Update-DatePlaceholder $master "9/22/2024"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i) "9/22/2024"
}

# ---------------------------------------------------------------------
# 2) Slide 4 - "TextBox 7" feature summary (trim trailing clause, and
#    shrink the auto-fit textbox to match the shorter text's height)
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tb7 = Get-ShapeByName $slide4 "TextBox 7"
$tb7.TextFrame.TextRange.Text = "Our application boasts several key features such as secure registration, meal logging, expiration alerts and sustainability tips."
$tb7.Height = 1754326 / 12700

# ---------------------------------------------------------------------
# 3) Slide 8 - "TextBox 9" Diet Tracker summary (trim trailing clause,
#    and shrink the auto-fit textbox to match the shorter text's height)
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$tb9 = Get-ShapeByName $slide8 "TextBox 9"
$tb9.TextFrame.TextRange.Text = "By utilizing Diet Tracker, users can simplify their meal tracking, adopt healthier eating habits, minimize food waste, and cultivate a mindset of sustainability."
$tb9.Height = 2308324 / 12700
